$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Shopenzer Testcases")
$ws2 = $wb.Worksheets.Item("Testscearnios")

# Update the Team ID
$ws1.Range("F2").Value = "PNT2022TMID53380"

# Update the team member names
$ws1.Range("N6").Value = "Ritunjay M"
$ws1.Range("N7").Value = "Praveen Raagul R"
$ws1.Range("N8").Value = "Pradeep V"
$ws1.Range("N9").Value = "Munish Kumar S"

# Switch the active sheet/selection to the "Shopenzer Testcases" sheet
$ws1.Activate()
$ws1.Range("N9").Select()
